# "Dummy First Record in All Excel Files"
# Adds a second (dummy/placeholder) data row to the CatalogProductUpload
# sample sheet, formats it with a yellow highlight band, wraps text on a
# couple of columns, applies a TRUE/FALSE custom number format to the two
# boolean "feature flag" columns, and widens several header columns so the
# new sample values are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Widen columns so the new dummy data is readable
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 19.9    # Name
$ws.Columns("B").ColumnWidth = 15.63   # ParentCode
$ws.Columns("C").ColumnWidth = 15.83   # Sku
$ws.Columns("N").ColumnWidth = 12.33   # Brand
$ws.Range($ws.Columns("R"), $ws.Columns("T")).ColumnWidth = 20.4   # Id / SellerId / CategoryId

# ---------------------------------------------------------------------
# 2) Write the dummy product record into row 2
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Dummy Product Name "
$ws.Range("B2").Value = "LR_AAA_12345"
$ws.Range("C2").Value = "LR_AAA_12345"
$ws.Range("D2").Formula = "=FALSE()"
$ws.Range("E2").Formula = "=FALSE()"
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "FIXED"
$ws.Range("K2").Value = 0
$ws.Range("L2").Formula = "=FALSE()"
# M2 (ModelNo) intentionally left blank for the dummy record
$ws.Range("N2").Value = "Dummy Brand"
$ws.Range("O2").Value = "red"
$ws.Range("P2").Value = "L"
# Q2 (Images) intentionally left blank for the dummy record
$ws.Range("R2").Value = "1a2b3c4d5e-1234-1234"
$ws.Range("S2").Value = "1a2b3c4d5e-1234-1234"
$ws.Range("T2").Value = "1a2b3c4d5e-1234-1234"
$ws.Range("U2").Formula = "=TRUE()"
$ws.Range("V2").Value = 1
$ws.Range("W2").Formula = "=FALSE()"
$ws.Range("X2").Value = 0

# ---------------------------------------------------------------------
# 3) Format row 2: 11pt Calibri, yellow fill across the whole row,
#    wrapped text for the two code columns, and a TRUE/FALSE display
#    format on the two boolean "feature" flags.
# ---------------------------------------------------------------------
$black = 0          # RGB(0,0,0)
$yellow = 62207     # RGB(255,242,0) == BGR 0x00F2FF used by the workbook's fill

$dataRow = $ws.Range("A2:X2")
$dataRow.Font.Size = 11
$dataRow.Font.Name = "Calibri"
$dataRow.Font.Color = $black
$dataRow.Interior.Color = $yellow
$ws.Rows(2).RowHeight = 15

$ws.Range("B2:C2").WrapText = $true

$ws.Range("U2").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("W2").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Apply the (no-fill) 11pt Calibri font as the row's baseline style too,
# so any untouched cell in the row still picks up the new default look.
$ws.Rows(2).Font.Size = 11
$ws.Rows(2).Font.Name = "Calibri"
$ws.Rows(2).Font.Color = $black

# ---------------------------------------------------------------------
# 4) Update the view: scroll/selection moves to the new Brand cell (N2)
# ---------------------------------------------------------------------
$ws.Range("N2").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

Write-Output "Dummy record written to row 2"
